$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("edittedData")

$headers = @(
    "ID",
    "MUID",
    "TERM",
    "COMPANY_ID",
    "ACTIVITY",
    "SALARY",
    "CITY",
    "STATE",
    "COUNTRY",
    "REGID",
    "WORK_REG",
    "WORK_GRADE",
    "GRADING_REG",
    "GRADING_GRADE",
    "EMPLOYER_EVAL_DATE",
    "EMPLOYER_EVAL",
    "EMPLOYER_AUTH",
    "STUDENT_EVAL_DATE",
    "STUDENT_EVAL",
    "STUDENT_EVAL_DATE"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# The source row was originally filled across A1:U1, but the last cell (U1)
# never received a value, leaving it blank while still being part of the
# used range. Touch and revert its formatting so the sheet's dimension
# extends through column U without adding any content there.
$ws.Cells.Item(1, 21).Font.Bold = $true
$ws.Cells.Item(1, 21).Font.Bold = $false
